# The workbook's single table needs a new leading column inserted (shifting
# the existing "Category"/"Group"/"Treatment"/"Control" columns one slot to
# the right), the "n=" / "P=" expressions get a space added around the "="
# sign, and the newly vacated column A gets the bold/bordered/centered header
# formatting applied down through row 15 (matching the styling already used
# by the header row and by column A in the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every existing column one place to the right by inserting a new
# blank column A. Excel automatically moves all cell contents/styles/refs.
$ws.Columns("A:A").Insert()

# Fix up the header row text (spaces added around "=")
$ws.Range("D1").Value = "Treatment at T1 (n = 5080)"
$ws.Range("E1").Value = "Control at T2 (n = 745)"

# Fix up the category label text (spaces added around "=")
$ws.Range("B3").Value = "Gender (P = 0.006)"
$ws.Range("B10").Value = "Interested in News (P = 0.000)"

# Apply the bold/centered/bordered header-style formatting (as used by row 1
# and originally by column A) down column A for rows 2-15, leaving the
# cells themselves empty.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A2:A15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

Write-Host "Applied column insert and header text fixes"
